$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'basketball leg sleeve youth padded'
$ws.Range('A2').Value = 'knee pad gym'
$ws.Range('A3').Value = 'work need pads'
$ws.Range('A4').Value = 'softball material'
$ws.Range('A5').Value = 'compression spandex men'
$ws.Range('A6').Value = 'football leggings men'
$ws.Range('A7').Value = 'knees bees'
$ws.Range('A8').Value = 'snowboarding pants youth'
$ws.Range('A9').Value = 'knee length pants'
$ws.Range('A10').Value = 'medias de basketball'
$ws.Range('A11').Value = 'padded sliding shorts youth'
$ws.Range('A12').Value = 'men gym tight pants'
$ws.Range('A13').Value = 'compression tight leggings'
$ws.Range('A14').Value = 'hockey leggings girls'
$ws.Range('A15').Value = 'boys black tight pants'
$ws.Range('A16').Value = 'yoga pad'
$ws.Range('A17').Value = 'leg compression pants men'
$ws.Range('A18').Value = 'wrestling shorts men'
$ws.Range('A19').Value = 'yoga pants men big and tall'
$ws.Range('A20').Value = 'wrestling knee'
$ws.Range('A21').Value = 'youth sliding shorts'
$ws.Range('A22').Value = 'youth hockey pants'
$ws.Range('A23').Value = 'sliding shorts baseball'
$ws.Range('A24').Value = 'knee pads replacement'
$ws.Range('A25').Value = 'thread protector 1/2 x 28'
$ws.Range('A26').Value = 'running tights youth'
$ws.Range('A27').Value = 'mesh leggings men'
$ws.Range('A28').Value = 'boys hockey pants'
$ws.Range('A29').Value = 'baseball pants youth boys black'
$ws.Range('A30').Value = 'basketball tight shorts for boys'
$ws.Range('A31').Value = 'sports compression leggings'
$ws.Range('A32').Value = 'knee pads for yoga'
$ws.Range('A33').Value = 'padded work pants mens'
$ws.Range('A34').Value = 'baseball pants for boys'
$ws.Range('A35').Value = 'youth spandex'
$ws.Range('A36').Value = 'women knee pads for work'
$ws.Range('A37').Value = 'big and tall mens compression pants'
$ws.Range('A38').Value = 'padded compression shorts youth'
$ws.Range('A39').Value = 'boys tights youth'
$ws.Range('A40').Value = 'men sheer pants'
$ws.Range('A41').Value = 'football girdle with pads for men'
$ws.Range('A42').Value = 'knee pads for work men'
$ws.Range('A43').Value = 'yoga position chart'
$ws.Range('A44').Value = 'youth football leggings boys'
$ws.Range('A45').Value = 'knee compression sleeve pad'
$ws.Range('A46').Value = 'black baseball pants mens'
$ws.Range('A47').Value = 'mens leggings shorts'
$ws.Range('A48').Value = 'youth baseball pants knee high'
$ws.Range('A49').Value = 'elastic waist baseball pants'
$ws.Range('A50').Value = 'adult tights'
$ws.Range('A51').Value = 'knee pads for work black'
$ws.Range('A52').Value = 'arthritis test'
$ws.Range('A53').Value = 'cycling sweat guard'
$ws.Range('A54').Value = 'volleyball knee pads extra large'
$ws.Range('A55').Value = 'size 5 basketball'
$ws.Range('A56').Value = 'athletic knee compression'
$ws.Range('A57').Value = 'knee pads for man'
$ws.Range('A58').Value = 'youth basketball knee sleeve'
$ws.Range('A59').Value = 'knee pads for mountain biking'
$ws.Range('A60').Value = 'best knee pads for work'
$ws.Range('A61').Value = 'hip pads for men'
$ws.Range('A62').Value = 'girl compression pants'
$ws.Range('A63').Value = 'compression shorts basketball'
$ws.Range('A64').Value = 'basketball knee sleeve boys'
$ws.Range('A65').Value = 'mens pants big and tall'
$ws.Range('A66').Value = 'below knee shorts men'
$ws.Range('A67').Value = 'youth wrestling shorts'
$ws.Range('A68').Value = 'capri spandex'
$ws.Range('A69').Value = 'yoga pads for hands'
$ws.Range('A70').Value = 'football pants adult with pads'
$ws.Range('A71').Value = 'boys youth compression pants'
$ws.Range('A72').Value = 'exercise kneeling pad'
$ws.Range('A73').Value = 'knee pad volleyball'
$ws.Range('A74').Value = 'knee pads working'
$ws.Range('A75').Value = 'baseball softball pants'
$ws.Range('A76').Value = 'boys knee pads volleyball'
$ws.Range('A77').Value = 'sliding workout pads'
$ws.Range('A78').Value = 'knee pads for youth'
$ws.Range('A79').Value = 'black legging for men'
$ws.Range('A80').Value = 'cheap leggings for men'
$ws.Range('A81').Value = 'little boys compression leggings'
$ws.Range('A82').Value = 'basketball compression gear'
$ws.Range('A83').Value = 'men sport pants'
$ws.Range('A84').Value = 'cold knee pad'
$ws.Range('A85').Value = 'black compression shorts for men'
$ws.Range('A86').Value = 'calf sleeves for men basketball'
$ws.Range('A87').Value = 'mens work knee pads'
$ws.Range('A88').Value = 'the bees knees'
$ws.Range('A89').Value = 'athletic leggings for men'
$ws.Range('A90').Value = 'compression pants for boys'
$ws.Range('A91').Value = 'patella knee pads'
$ws.Range('A92').Value = 'lacrosse tights'
$ws.Range('A93').Value = 'boys leggings youth'
$ws.Range('A94').Value = 'boy sport tights'
$ws.Range('A95').Value = 'volleyball pants'
$ws.Range('A96').Value = 'knee pads for wrestling'
$ws.Range('A97').Value = 'football padded shorts for men'
$ws.Range('A98').Value = 'mens basketball outdoor'
$ws.Range('A99').Value = 'compression pants size'
$ws.Range('A100').Value = 'mens knee pads construction'
